$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 76: tide pool samples 19 and 20 from OR (second entry added)
$ws.Range("A76").Value = 43805
$ws.Range("B76").Value = 2206.5185511325199
$ws.Range("C76").Value = 2207.0300000000002
$ws.Range("D76").Formula = "=100*(B76-C76)/C76"
$ws.Range("E76").Value = 169
$ws.Range("F76").Value = "Crm opened 11/19/2020"

# Match the date format used by the rest of column A (copy style from A75)
$ws.Range("A75").Copy() | Out-Null
$ws.Range("A76").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the view to match the new selection / scroll position
$ws.Range("A58").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F75:F76").Select() | Out-Null
